# LogExcel.xlsx — "Add files via upload" edit
#
# 1. Fill in the two previously-empty cells in row 4 (B4, C4) with their
#    new numeric values.
# 2. Update the sheet's view state: the window is scrolled down so row 4
#    is at the top, and the active selection moves from C2 to C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: row 4 gains values in columns B and C ---------------------
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 2

# --- View-state edit: scroll position + selection --------------------------
$win = $excel.ActiveWindow

# Scroll the window so A4 becomes the top-left visible cell.
$win.ScrollRow = 4
$win.ScrollColumn = 1

# Move the active selection to C23.
$ws.Range("C23").Select()
